$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.05662722568076514
$ws.Range("C2").Value = 2.181838915208745
$ws.Range("D2").Value = 17.88885873886894
$ws.Range("E2").Value = 4.229522282583335
$ws.Range("F2").Value = 4.333582201942138

# Row 3 (Q0)
$ws.Range("B3").Value = 0.927092264295462
$ws.Range("C3").Value = 2.285903137075404
$ws.Range("D3").Value = 20.1805982684956
$ws.Range("E3").Value = 4.492282078019545
$ws.Range("F3").Value = 4.499016614473818

# Row 4 (Q1)
$ws.Range("B4").Value = -0.2082090657517003
$ws.Range("C4").Value = 1.005865323363158
$ws.Range("D4").Value = 3.987134668739914
$ws.Range("E4").Value = 1.996781076818366
$ws.Range("F4").Value = 2.040314374195741
$ws.Range("G4").Value = 19
